$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 510, pushing the existing rows
# 510-534 down to 512-536 (weekly data roll: two new sampling dates
# added at the front of this date-grouped block).
$ws.Range("A510:R511").Insert()

# Row 510 - new record: Crespo record, date 44509
$ws.Range("A510").Value = 10
$ws.Range("B510").Value = "Vega Modelo de Temuco"
$ws.Range("C510").Value = "La Araucanía"
$ws.Range("D510").Value = 44509
$ws.Range("E510").Value = 9
$ws.Range("F510").Value = 100112006
$ws.Range("G510").Value = "Repollo"
$ws.Range("H510").Value = "Crespo record"
$ws.Range("I510").Value = "Primera"
$ws.Range("J510").Value = 400
$ws.Range("K510").Value = 1000
$ws.Range("L510").Value = 1000
$ws.Range("M510").Value = 1000
$ws.Range("N510").Value = "`$/unidad"
$ws.Range("O510").Value = "Región Metropolitana"
$ws.Range("P510").Value = 1000
$ws.Range("Q510").Value = 1
$ws.Range("R510").Value = "Hortaliza"

# Row 511 - new record: Morada(o), date 44509
$ws.Range("A511").Value = 10
$ws.Range("B511").Value = "Vega Modelo de Temuco"
$ws.Range("C511").Value = "La Araucanía"
$ws.Range("D511").Value = 44509
$ws.Range("E511").Value = 9
$ws.Range("F511").Value = 100112006
$ws.Range("G511").Value = "Repollo"
$ws.Range("H511").Value = "Morada(o)"
$ws.Range("I511").Value = "Primera"
$ws.Range("J511").Value = 100
$ws.Range("K511").Value = 1000
$ws.Range("L511").Value = 1200
$ws.Range("M511").Value = 1100
$ws.Range("N511").Value = "`$/unidad"
$ws.Range("O511").Value = "Región Metropolitana"
$ws.Range("P511").Value = 1100
$ws.Range("Q511").Value = 1
$ws.Range("R511").Value = "Hortaliza"

# Ensure the date columns (D510, D511) carry the same date number format
# as the rest of the D column (style index 2 in styles.xml).
$ws.Range("D510:D511").NumberFormat = $ws.Range("D509").NumberFormat
